# "Loan RBI, Variable Instalments"
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, pushing the old N:P ("Late", heading, Outstanding) columns to
#   O:Q, and give the new column the same width as column M.
# - Make "Repayment schedule" the active sheet/tab, with K16 selected.
# - Leave the "Edit Repayment Schedule" sheet with B8 selected (and no
#   longer the active tab).

$wb = $excel.ActiveWorkbook

$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsEdit.Range("B8").Select()

$wsRepay = $wb.Worksheets.Item("Repayment schedule")

$colWidthM = $wsRepay.Columns("M:M").ColumnWidth
$wsRepay.Columns("N:N").Insert()
$wsRepay.Columns("N:N").ColumnWidth = $colWidthM

$wsRepay.Activate()
$wsRepay.Range("K16").Select()
